$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price & volume columns) as text, preserving original formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.581.48'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.124.69'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.24%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.26'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.63'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.95%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.118.40'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.30'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.51%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.04%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.28'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.642.87'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.25%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.683.39'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.128.40'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.46%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '474.78'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.53'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.710'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.84'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.34'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.84%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.45'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.120'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -7.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.01'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.43%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -10.23%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.67'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.16'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.36%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0787'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.44%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.83'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '455.83'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.99'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -8.18%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.49%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.30'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.848.07'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.29'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.65%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.57%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.20'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.84%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.88'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.46%  '
